$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: B7/C7 text updates (未登録 due to unclear inspection contents)
$ws.Range("B7").Value = '**大阪健康安全基盤研究所** <br> Ready-to-eat 食品製造環境のリステリア汚染状況と食品の特性に応じた汚染制御についての考察 <br> （[大阪健康安全基盤研究所研究年報, p.107, 2020](http://www.iph.osaka.jp/s004/060/reserch_annual_report-7_2023.pdf) / 月刊 HACCP, 10, 20-25 (2022)）<br> 検査内容不明のため未登録'
$ws.Range("C7").Value = '未登録'

# Insert a new row at 54 (shifts old rows 54-64 down to 55-65),
# then populate the newly inserted row with the 広島市衛生研究所 entry.
$ws.Rows.Item(54).EntireRow.Insert()

# A54 is a purely-numeric-looking label ("2004"); assigning a plain string
# via .Value lets Excel's type-inference store it as a *number*, unlike the
# original file where every A-column year/range label is text. Write it as
# a text-producing formula, then collapse the formula to a static value via
# copy / paste-values (xlPasteValues = -4163) so it lands as a plain text
# cell with no residual formula and no extra number-format/style.
$aCell = $ws.Range("A54")
$aCell.Formula = '="2004"'
$aCell.Copy()
$aCell.PasteSpecial(-4163)

$ws.Range("B54").Value = '**広島市衛生研究所** <br> [鶏肉のカンピロバクター培養検査法の検討 -鶏肉の検査方法別検出感度および検出率の比較]() <br>（広島市衛研年報 24, 49-54, 2005）'
$ws.Range("C54").Value = '未登録'

